# feat(matriz-adjacencia): cria matriz de adjacencia para grafo dirigido e
# nao ponderado, e para grafo nao dirigido e nao ponderado
#
# - sheet "vertice": drop vertices E and F, keep only A, B, C, D
# - sheet "aresta": replace the 5 directed/weighted edges with the 6
#   edges used by the new (unweighted) adjacency matrices, all with
#   weight 0: AA(A,B) AB(A,B) BC(B,C) BD(B,D) CB(C,B) DC(D,C)

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("vertice")
$ws2 = $wb.Worksheets.Item("aresta")

# --- sheet "vertice": remove rows 5 and 6 (vertices E, F) ---
$ws1.Range("A5:A6").EntireRow.Delete()

# --- sheet "aresta": overwrite the 5 existing rows and append a 6th ---
$ws2.Range("A1").Value = "AA"
$ws2.Range("B1").Value = "A"
$ws2.Range("C1").Value = "B"
$ws2.Range("D1").Value = 0

$ws2.Range("A2").Value = "AB"
$ws2.Range("B2").Value = "A"
$ws2.Range("C2").Value = "B"
$ws2.Range("D2").Value = 0

$ws2.Range("A3").Value = "BC"
$ws2.Range("B3").Value = "B"
$ws2.Range("C3").Value = "C"
$ws2.Range("D3").Value = 0

$ws2.Range("A4").Value = "BD"
$ws2.Range("B4").Value = "B"
$ws2.Range("C4").Value = "D"
$ws2.Range("D4").Value = 0

$ws2.Range("A5").Value = "CB"
$ws2.Range("B5").Value = "C"
$ws2.Range("C5").Value = "B"
$ws2.Range("D5").Value = 0

$ws2.Range("A6").Value = "DC"
$ws2.Range("B6").Value = "D"
$ws2.Range("C6").Value = "C"
$ws2.Range("D6").Value = 0

# --- selections: "vertice" keeps a non-active selection at B5, while
#     "aresta" stays the active/tabSelected sheet with selection C6 ---
[void]$ws1.Range("B5").Select()
[void]$ws2.Select()
[void]$ws2.Range("C6").Select()
